# Append four new paragraphs to the end of the document, after the
# existing trailing empty paragraph:
#   "I am adding a new text here."
#   (empty)
#   (empty)
#   "Something new."
#
# Word's Range.InsertXML merges the *last* <w:p> fragment's paragraph
# mark with the target (collapsed) range's own paragraph mark - exactly
# like typing text followed by Enter does. A leading empty <w:p/>
# fragment is included so that merge lands on a throwaway paragraph
# mark, leaving the document's existing trailing paragraph untouched
# (still empty) and producing genuinely empty <w:p/> paragraphs (no
# stray runs) for the two blank lines, matching a direct XML edit.

$d = $word.ActiveDocument

$rng = $d.Content
$rng.Collapse(0)

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$xml = "<w:p $w/>" +
       "<w:p $w><w:r><w:t>I am adding a new text here.</w:t></w:r></w:p>" +
       "<w:p $w/>" +
       "<w:p $w/>" +
       "<w:p $w><w:r><w:t>Something new.</w:t></w:r></w:p>"

$rng.InsertXML($xml)

Write-Host "Paragraph count now:" $d.Paragraphs.Count
